$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices, volume %, and swapped row
# order for Stacks/Mantle and OKB/VeChain) as produced by the scheduled
# GitHub Actions data refresh. Price/volume cells are stored as text in
# the workbook, so force text format before writing values that would
# otherwise be auto-converted to numbers by Excel (e.g. "0.999", "1.20").

$ws.Range("D2").Value = '67.444.47'
$ws.Range("E2").Value = '  -2.95%  '
$ws.Range("D3").Value = '3.259.38'
$ws.Range("E3").Value = '  -5.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.69'
$ws.Range("E5").Value = '  -3.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.14'
$ws.Range("E6").Value = '  -10.54%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '3.251.24'
$ws.Range("E8").Value = '  -5.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  -8.49%  '
$ws.Range("E10").Value = '  -10.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.69'
$ws.Range("E11").Value = '  -4.95%  '
$ws.Range("E12").Value = '  -10.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  -8.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.26'
$ws.Range("E14").Value = '  -13.75%  '
$ws.Range("D15").Value = '3.778.63'
$ws.Range("E15").Value = '  -5.86%  '
$ws.Range("D16").Value = '67.435.26'
$ws.Range("D17").Value = '3.257.45'
$ws.Range("E17").Value = '  -5.60%  '
$ws.Range("E18").Value = '  -5.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '529.65'
$ws.Range("E19").Value = '  -8.55%  '
$ws.Range("E20").Value = '  -13.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.94'
$ws.Range("E21").Value = '  -12.98%  '
$ws.Range("E22").Value = '  -10.90%  '
$ws.Range("E23").Value = '  -11.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.48'
$ws.Range("E24").Value = '  -10.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.48'
$ws.Range("E25").Value = '  -11.30%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.21'
$ws.Range("E27").Value = '  -11.57%  '
$ws.Range("E28").Value = '  -11.92%  '
$ws.Range("E29").Value = '  -7.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.06'
$ws.Range("E30").Value = '  -11.42%  '
$ws.Range("B31").Value = 'Mantle'
$ws.Range("C31").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.20'
$ws.Range("E31").Value = '  -3.86%  '
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.67'
$ws.Range("E32").Value = '  -5.38%  '
$ws.Range("E33").Value = '  -15.98%  '
$ws.Range("E34").Value = '  -13.29%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '510.52'
$ws.Range("E36").Value = '  -11.96%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.15'
$ws.Range("E37").Value = '  -5.18%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0438'
$ws.Range("E38").Value = '  -7.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0853'
$ws.Range("E39").Value = '  -10.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.93'
$ws.Range("E40").Value = '  -15.39%  '
$ws.Range("E41").Value = '  -10.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.75'
$ws.Range("E42").Value = '  -13.00%  '
$ws.Range("D43").Value = '2.929.80'
$ws.Range("E43").Value = '  -9.94%  '
$ws.Range("E44").Value = '  -10.25%  '
$ws.Range("D45").Value = '0.0₃0587'
$ws.Range("E45").Value = '  -15.22%  '
$ws.Range("E46").Value = '  -8.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.65'
$ws.Range("E47").Value = '  -14.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.31'
$ws.Range("E49").Value = '  -17.26%  '
$ws.Range("E50").Value = '  -10.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '123.74'
